$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the old "_GoBack" bookmark that sits between " on " and
#    "fit" (the diff deletes these two tags outright).
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2) Re-create a "_GoBack" bookmark spanning the paragraph that holds
#    the "ridershipByHour.png" picture (between "Other Visualizations:"
#    and the "Plot:  Ridership by hour, via resampling" paragraph).
#    This reproduces the net effect of the diff, which re-homes the
#    bookmark there and shifts the other bookmark ids down by one.
# ------------------------------------------------------------------
$targetPara = $null
$idx = 0
foreach ($p in $d.Paragraphs) {
    $idx = $idx + 1
    $t = $p.Range.Text
    if ($t -like "*Other Visualizations:*") {
        $targetPara = $p
    }
}

# $targetPara is "Other Visualizations:" -- the very next paragraph is
# the one that contains the inline picture.
$found = $false
$idx = 0
$pictureRange = $null
foreach ($p in $d.Paragraphs) {
    $idx = $idx + 1
    if ($found) {
        $pictureRange = $p.Range
        $found = $false
    }
    if ($p.Range.Start -eq $targetPara.Range.Start) {
        $found = $true
    }
}

$d.Bookmarks.Add("_GoBack", $pictureRange)

# ------------------------------------------------------------------
# 3) Fix the typo: ".  The low p-value to" -> ".  The low p-value led to"
# ------------------------------------------------------------------
$d.Content.Find.Execute(".  The low p-value to", $true, $false, $false, $false, $false, $true, 1, $false, ".  The low p-value led to", 2)
